# Generate Report for Handoff
# Update the localization status report: mark rows 4-7 (the e2e "ht" set) as
# newly handed off by bumping their Priority to "ht" and refreshing the
# "Latest Handoff Datetime" timestamp, for both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

for ($row = 4; $row -le 7; $row++) {
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = "2016-08-13 04:36:13"

    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = "2016-08-13 04:36:21"
}
